$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 23550
$ws.Range("J21").Value = 23550
$ws.Range("L21").Value = 23550
$ws.Range("N21").Value = -24486

$ws.Range("H23").Value = 23550
$ws.Range("J23").Value = 23550
$ws.Range("L23").Value = 23550
$ws.Range("N23").Value = -24018

$ws.Range("H32").Value = 858.125
$ws.Range("I32").Value = 975.55554
$ws.Range("J32").Value = 787.6667
$ws.Range("K32").Value = 975.55554
$ws.Range("L32").Value = 787.6667
$ws.Range("M32").Value = -649.55554
$ws.Range("N32").Value = -1439.6667

$ws.Range("H34").Value = 2520
$ws.Range("I34").Value = 2520
$ws.Range("K34").Value = 2520
$ws.Range("M34").Value = -2317

$ws.Range("H36").Value = 2520
$ws.Range("I36").Value = 2520
$ws.Range("K36").Value = 2520
$ws.Range("M36").Value = -1805

$ws.Range("H52").Value = 1490
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1490
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").Value = 4470
$ws.Range("N52").Value = -4790

$ws.Range("H132").Value = 2520.2131
$ws.Range("I132").Value = 1998.8246
$ws.Range("J132").Value = 9950
$ws.Range("K132").Value = 5996.4738
$ws.Range("L132").Value = 29850
$ws.Range("M132").Value = -3466.4738
$ws.Range("N132").Value = -34910

$ws.Range("H137").Value = 5129061
$ws.Range("I137").Value = 823.875
$ws.Range("J137").Value = 13334241
$ws.Range("K137").Value = 2471.625
$ws.Range("L137").Value = 40002723
$ws.Range("M137").Value = 78.375
$ws.Range("N137").Value = -40007823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3393.7144
$ws.Range("I45").Value = 6006
$ws.Range("K45").Value = 6006
$ws.Range("M45").Value = -5629

$ws.Range("H61").Value = 7463799.5
$ws.Range("I61").Value = 9435003
$ws.Range("J61").Value = 1388.1428
$ws.Range("K61").Value = 9435003
$ws.Range("L61").Value = 1388.1428
$ws.Range("M61").Value = -9434791
$ws.Range("N61").Value = -1812.1428

$ws.Range("H132").Value = 6251810.5
$ws.Range("I132").Value = 8930051
$ws.Range("J132").Value = 2582.6667
$ws.Range("K132").Value = 26790153
$ws.Range("L132").Value = 7748.000100000001
$ws.Range("M132").Value = -26787623
$ws.Range("N132").Value = -12808.0001

$ws.Range("H136").Value = 7463799.5
$ws.Range("I136").Value = 9435003
$ws.Range("J136").Value = 1388.1428
$ws.Range("K136").Value = 28305009
$ws.Range("L136").Value = 4164.428400000001
$ws.Range("M136").Value = -28302459
$ws.Range("N136").Value = -9264.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2034.6522
$ws.Range("I134").Value = 1282.375
$ws.Range("K134").Value = 3847.125
$ws.Range("M134").Value = -1312.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 991.4211
$ws.Range("I16").Value = 930.0833
$ws.Range("K16").Value = 930.0833
$ws.Range("M16").Value = -643.0833

$ws.Range("H31").Value = 5558900
$ws.Range("I31").Value = 3950.524
$ws.Range("J31").Value = 18520448
$ws.Range("K31").Value = 3950.524
$ws.Range("L31").Value = 18520448
$ws.Range("M31").Value = -3655.524
$ws.Range("N31").Value = -18521038

$ws.Range("H34").Value = 5558900
$ws.Range("I34").Value = 3950.524
$ws.Range("J34").Value = 18520448
$ws.Range("K34").Value = 3950.524
$ws.Range("L34").Value = 18520448
$ws.Range("M34").Value = -3748.524
$ws.Range("N34").Value = -18520852

$ws.Range("H113").Value = 991.4211
$ws.Range("I113").Value = 930.0833
$ws.Range("K113").Value = 930.0833
$ws.Range("M113").Value = 1239.9167

$ws.Range("H132").Value = 22730206
$ws.Range("I132").Value = 31252314
$ws.Range("J132").Value = 4585.3335
$ws.Range("K132").Value = 93756942
$ws.Range("L132").Value = 13756.0005
$ws.Range("M132").Value = -93754412
$ws.Range("N132").Value = -18816.0005

$ws.Range("H140").Value = 45459.832
$ws.Range("J140").Value = 45459.832
$ws.Range("L140").Value = 45459.832
$ws.Range("N140").Value = -55819.832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 136.54546
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 136.54546
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 409.63638
$ws.Range("N23").Value = -879.6363799999999

$ws.Range("H68").Value = 1023.7778
$ws.Range("I68").Value = 846
$ws.Range("J68").Value = 1074.5714
$ws.Range("K68").Value = 2538
$ws.Range("L68").Value = 3223.7142
$ws.Range("M68").Value = -1727
$ws.Range("N68").Value = -4845.7142

$ws.Range("H71").Value = 1023.7778
$ws.Range("I71").Value = 846
$ws.Range("J71").Value = 1074.5714
$ws.Range("K71").Value = 7614
$ws.Range("L71").Value = 9671.142600000001
$ws.Range("M71").Value = -3558
$ws.Range("N71").Value = -17783.1426

$ws.Range("H92").Value = 312.5
$ws.Range("I92").Value = 293.33334
$ws.Range("J92").Value = 370
$ws.Range("K92").Value = 880.0000200000001
$ws.Range("L92").Value = 1110
$ws.Range("M92").Value = 367.9999799999999
$ws.Range("N92").Value = -3606

$ws.Range("H122").Value = 1948.1111
$ws.Range("J122").Value = 1599.75
$ws.Range("L122").Value = 14397.75
$ws.Range("N122").Value = -19297.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -21064

$ws.Range("H70").Value = 6600
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 8333.333000000001
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 8333.333000000001
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -8873.333000000001

$ws.Range("H73").Value = 6600
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 8333.333000000001
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 8333.333000000001
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -10205.333

$ws.Range("H113").Value = 63699.688
$ws.Range("I113").Value = 91745.55
$ws.Range("J113").Value = 1998.8
$ws.Range("K113").Value = 91745.55
$ws.Range("L113").Value = 1998.8
$ws.Range("M113").Value = -89575.55
$ws.Range("N113").Value = -6338.8

$ws.Range("H132").Value = 2735.9104
$ws.Range("I132").Value = 1898.3846
$ws.Range("K132").Value = 5695.1538
$ws.Range("M132").Value = -3165.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 23669.143
$ws.Range("J50").Value = 23669.143
$ws.Range("L50").Value = 23669.143
$ws.Range("N50").Value = -24943.143

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H54").Value = 34870.25
$ws.Range("J54").Value = 34870.25
$ws.Range("L54").Value = 34870.25
$ws.Range("N54").Value = -36158.25

$ws.Range("H55").Value = 1100
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1100
$ws.Range("K55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("M55").Value = 1100
$ws.Range("N55").Value = -1446

$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622

$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112

$ws.Range("H132").Value = 8172.5854
$ws.Range("I132").Value = 4967.95
$ws.Range("J132").Value = 11224.619
$ws.Range("K132").Value = 14903.85
$ws.Range("L132").Value = 33673.857
$ws.Range("M132").Value = -12373.85
$ws.Range("N132").Value = -38733.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 2000
$ws.Range("I43").Value = 2000
$ws.Range("K43").Value = 2000
$ws.Range("M43").Value = -1851
